# This workbook ("Lich_Profits") caches market-board snapshot values
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) computed by an
# external scheduled job. This run refreshes those cached numbers in place
# -- there are no formulas on these sheets, so each changed cell is just a
# literal numeric overwrite (plus one cell that drops out of HQ pricing and
# one that newly gains an NQ profit figure).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 2139.6667
$ws.Range("J19").Value = 2578.3076
$ws.Range("L19").Value = 2578.3076
$ws.Range("N19").Value = -2928.3076
# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 1733.7693
$ws.Range("I100").Value = 1621.7273
$ws.Range("J100").Value = 2350
$ws.Range("K100").Value = 1621.7273
$ws.Range("L100").Value = 2350
$ws.Range("M100").Value = -1080.7273
$ws.Range("N100").Value = -3432
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2884.182
$ws.Range("I132").Value = 2883.2856
$ws.Range("J132").Value = 2887.6667
$ws.Range("K132").Value = 8649.856800000001
$ws.Range("L132").Value = 8663.000100000001
$ws.Range("M132").Value = -6119.856800000001
$ws.Range("N132").Value = -13723.0001
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 912.7
$ws.Range("I135").Value = 900.2941
$ws.Range("J135").Value = 983
$ws.Range("K135").Value = 8102.6469
$ws.Range("L135").Value = 8847
$ws.Range("M135").Value = -5567.6469
$ws.Range("N135").Value = -13917
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2249437
$ws.Range("I137").Value = 3612654.8
$ws.Range("J137").Value = 4137
$ws.Range("K137").Value = 10837964.4
$ws.Range("L137").Value = 12411
$ws.Range("M137").Value = -10835414.4
$ws.Range("N137").Value = -17511

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 62793.03
$ws.Range("I74").Value = 66646.07000000001
$ws.Range("K74").Value = 66646.07000000001
$ws.Range("M74").Value = -65772.07000000001
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 62793.03
$ws.Range("I77").Value = 66646.07000000001
$ws.Range("K77").Value = 333230.35
$ws.Range("M77").Value = -328862.35

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1992.64
$ws.Range("I107").Value = 2032.6818
$ws.Range("K107").Value = 2032.6818
$ws.Range("M107").Value = -112.6818000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 3274.7742
$ws.Range("I58").Value = 2913.48
$ws.Range("K58").Value = 2913.48
$ws.Range("M58").Value = -2710.48
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 3439.8823
$ws.Range("J62").Value = 3537.4
$ws.Range("L62").Value = 3537.4
$ws.Range("N62").Value = -4785.4
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 3439.8823
$ws.Range("J65").Value = 3537.4
$ws.Range("L65").Value = 17687
$ws.Range("N65").Value = -23927
# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 3457.8125
$ws.Range("I122").Value = 3440.08
$ws.Range("K122").Value = 10320.24
$ws.Range("M122").Value = -7870.24
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1790.4736
$ws.Range("I132").Value = 1702.2858
$ws.Range("J132").Value = 2037.4
$ws.Range("K132").Value = 5106.857400000001
$ws.Range("L132").Value = 6112.200000000001
$ws.Range("M132").Value = -2576.857400000001
$ws.Range("N132").Value = -11172.2
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 2802
$ws.Range("I134").Value = 2782.1333
$ws.Range("K134").Value = 8346.3999
$ws.Range("M134").Value = -5811.3999
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 3274.7742
$ws.Range("I136").Value = 2913.48
$ws.Range("K136").Value = 8740.440000000001
$ws.Range("M136").Value = -6190.440000000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 913.3333
$ws.Range("I5").Value = 945.25
$ws.Range("J5").Value = 849.5
$ws.Range("K5").Value = 2835.75
$ws.Range("L5").Value = 2548.5
$ws.Range("M5").Value = -2723.75
$ws.Range("N5").Value = -2772.5
# Row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 249.38461
$ws.Range("J92").Value = 250.25
$ws.Range("L92").Value = 750.75
$ws.Range("N92").Value = -3246.75
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 502.125
$ws.Range("J113").Value = 584.6
$ws.Range("L113").Value = 1753.8
$ws.Range("N113").Value = -6093.8
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1407.2821
$ws.Range("I131").Value = 922.8570999999999
$ws.Range("J131").Value = 1513.25
$ws.Range("K131").Value = 2768.5713
$ws.Range("L131").Value = 4539.75
$ws.Range("M131").Value = 2271.4287
$ws.Range("N131").Value = -14619.75
# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 913.3333
$ws.Range("I135").Value = 945.25
$ws.Range("J135").Value = 849.5
$ws.Range("K135").Value = 8507.25
$ws.Range("L135").Value = 7645.5
$ws.Range("M135").Value = -5972.25
$ws.Range("N135").Value = -12715.5
# Row 141: Ocean Explosion / Acqua Pazza
$ws.Range("H141").Value = 3036.7273
$ws.Range("I141").Value = 2790.5
$ws.Range("K141").Value = 8371.5
$ws.Range("M141").Value = -3191.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 4077.2
$ws.Range("I113").Value = 3997
$ws.Range("J113").Value = 4130.6665
$ws.Range("K113").Value = 3997
$ws.Range("L113").Value = 4130.6665
$ws.Range("M113").Value = -1827
$ws.Range("N113").Value = -8470.666499999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = 45
$ws.Range("N22").Value = -940
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 250
$ws.Range("J27").Value = 350
$ws.Range("K27").Value = 250
$ws.Range("L27").Value = 350
$ws.Range("M27").Value = -143
$ws.Range("N27").Value = -564
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 2332.6428
$ws.Range("I46").Value = 2486.6
$ws.Range("J46").Value = 2155
$ws.Range("K46").Value = 2486.6
$ws.Range("L46").Value = 2155
$ws.Range("M46").Value = -2298.6
$ws.Range("N46").Value = -2531
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 2320.6667
$ws.Range("J61").Value = 1877.5
$ws.Range("L61").Value = 1877.5
$ws.Range("N61").Value = -2281.5
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 2380.75
$ws.Range("I100").Value = 2476.5217
$ws.Range("K100").Value = 2476.5217
$ws.Range("M100").Value = -1935.5217
# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 2320.6667
$ws.Range("J113").Value = 1877.5
$ws.Range("L113").Value = 1877.5
$ws.Range("N113").Value = -6217.5
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3740.9
$ws.Range("J122").Value = 4277
$ws.Range("L122").Value = 12831
$ws.Range("N122").Value = -17731
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3027.2144
$ws.Range("I132").Value = 2887.2
$ws.Range("J132").Value = 4194
$ws.Range("K132").Value = 8661.599999999999
$ws.Range("L132").Value = 12582
$ws.Range("M132").Value = -6131.599999999999
$ws.Range("N132").Value = -17642

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 6: Burn Me Up / Hempen Cowl
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 48: In over Your Head / Linen Cowl
$ws.Range("H48").Value = 120000
$ws.Range("I48").Value = 100000
$ws.Range("K48").Value = 100000
$ws.Range("M48").Value = -99431
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 2885.7144
$ws.Range("I122").Value = 2885.7144
$ws.Range("K122").Value = 8657.143199999999
$ws.Range("M122").Value = -6207.143199999999
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1863.3334
$ws.Range("I132").Value = 1226.6666
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3679.9998
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1149.9998
$ws.Range("N132").Value = -12560
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1001835.5
$ws.Range("I136").Value = 1001835.5
$ws.Range("K136").Value = 3005506.5
$ws.Range("M136").Value = -3002956.5

